$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.297.32'
$ws.Range("E2").Value = '  +3.07%  '

$ws.Range("D3").Value = '1.923.21'
$ws.Range("E3").Value = '  +2.68%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  -0.95%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.12'
$ws.Range("E5").Value = '  +1.49%  '

$ws.Range("E6").Value = '  -0.82%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4863'
$ws.Range("E7").Value = '  +1.48%  '

$ws.Range("E8").Value = '  +2.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07399'
$ws.Range("E9").Value = '  +0.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9429'
$ws.Range("E10").Value = '  +0.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.96'
$ws.Range("E11").Value = '  +1.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07818'
$ws.Range("E12").Value = '  -0.98%  '

$ws.Range("D13").Value = '1.942.32'
$ws.Range("E13").Value = '  +3.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.527'
$ws.Range("E14").Value = '  +1.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.667'
$ws.Range("E15").Value = '  +0.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.73'
$ws.Range("E16").Value = '  +0.99%  '

$ws.Range("E17").Value = '  -0.87%  '

$ws.Range("E18").Value = '  -0.32%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.007'
$ws.Range("E19").Value = '  -0.69%  '

$ws.Range("D20").Value = '28.310.28'
$ws.Range("E20").Value = '  +3.05%  '

$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.173'

$ws.Range("D23").Value = '2.182.74'
$ws.Range("E23").Value = '  +3.13%  '

$ws.Range("E24").Value = '  +2.51%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.40'
$ws.Range("E25").Value = '  +1.31%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.931'
$ws.Range("E26").Value = '  -1.68%  '

$ws.Range("E27").Value = '  +0.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.106'
$ws.Range("E28").Value = '  +4.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.71'
$ws.Range("E29").Value = '  +0.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.987'

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08928'
$ws.Range("E31").Value = '  -0.04%  '

$ws.Range("E32").Value = '  +1.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.254'
$ws.Range("E33").Value = '  +3.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7749'
$ws.Range("E34").Value = '  +3.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.712'
$ws.Range("E35").Value = '  +2.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.710'
$ws.Range("E36").Value = '  +0.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02053'
$ws.Range("E37").Value = '  -0.72%  '

$ws.Range("E38").Value = '  -1.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5572'
$ws.Range("E39").Value = '  +3.92%  '

$ws.Range("E40").Value = '  +0.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.020'
$ws.Range("E41").Value = '  +0.60%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.068'
$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1534'
$ws.Range("E43").Value = '  +0.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.497'
$ws.Range("E44").Value = '  +0.89%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.78'
$ws.Range("E45").Value = '  +1.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4885'
$ws.Range("E46").Value = '  +0.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '107.14'
$ws.Range("E47").Value = '  +3.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.667'
$ws.Range("E49").Value = '  +0.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '69.21'
$ws.Range("E50").Value = '  +3.20%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06141'
$ws.Range("E51").Value = '  +0.55%  '
